$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild header row: user import table now logs in via "nis" instead of
# "email", gains a "kelompok" column, and keeps "password" as the last
# column. Order of assignment matters for new shared-string indices, so
# write in the same order the strings were first introduced upstream:
# password, then kelompok, then nis (name/A1 is unchanged).
$ws.Range("D1").Value = "password"
$ws.Range("C1").Value = "kelompok"
$ws.Range("B1").Value = "nis"

# Move the active selection, matching the saved view state.
$ws.Range("B3").Select()
